$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header H1: teacher_ID -> teachers_ID
$ws.Range("H1").Value = "teachers_ID"

# Convert numeric teacher_ID values into bracketed list-of-IDs strings
$ws.Range("H2").Value = "[5, 9]"
$ws.Range("H3").Value = "[8]"
$ws.Range("H4").Value = "[9]"
$ws.Range("H5").Value = "[15]"
$ws.Range("H6").Value = "[45]"
$ws.Range("H7").Value = "[32]"
$ws.Range("H8").Value = "[6]"
$ws.Range("H9").Value = "[51]"
$ws.Range("H10").Value = "[46]"
$ws.Range("H11").Value = "[13]"
$ws.Range("H12").Value = "[23, 26]"
$ws.Range("H13").Value = "[24, 19]"
$ws.Range("H14").Value = "[27]"
$ws.Range("H15").Value = "[18]"
$ws.Range("H16").Value = "[14]"
$ws.Range("H17").Value = "[19]"
$ws.Range("H18").Value = "[20]"
$ws.Range("H19").Value = "[21]"
$ws.Range("H20").Value = "[2, 5]"
$ws.Range("H21").Value = "[23]"
$ws.Range("H22").Value = "[6]"
$ws.Range("H23").Value = "[8]"
$ws.Range("H24").Value = "[13]"
$ws.Range("H25").Value = "[16]"
$ws.Range("H26").Value = "[37]"
$ws.Range("H27").Value = "[48]"
$ws.Range("H28").Value = "[50]"
$ws.Range("H29").Value = "[10, 13]"
$ws.Range("H30").Value = "[14, 18]"
$ws.Range("H31").Value = "[16]"
$ws.Range("H32").Value = "[7]"
$ws.Range("H33").Value = "[26]"
$ws.Range("H34").Value = "[23]"
$ws.Range("H35").Value = "[29]"
$ws.Range("H36").Value = "[31, 33]"
$ws.Range("H37").Value = "[45, 41]"
$ws.Range("H38").Value = "[19]"
$ws.Range("H39").Value = "[3]"
$ws.Range("H40").Value = "[49]"
$ws.Range("H41").Value = "[51]"
$ws.Range("H42").Value = "[16]"
$ws.Range("H43").Value = "[24]"
$ws.Range("H44").Value = "[27]"
$ws.Range("H45").Value = "[29]"
$ws.Range("H46").Value = "[27, 32]"
$ws.Range("H47").Value = "[30, 33]"
$ws.Range("H48").Value = "[22]"
$ws.Range("H49").Value = "[34]"
$ws.Range("H50").Value = "[37]"
$ws.Range("H51").Value = "[44]"
$ws.Range("H52").Value = "[42]"

# Restore UI selection/scroll state to match the authored edit
$ws.Range("H52").Select()
